# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several crafting-job
# sheets in the Hades_Profits workbook. Pure data refresh, no structural
# changes -- only cell values on already-existing rows are touched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1007.8205
$ws.Range("J17").Value = 1007.8205
$ws.Range("L17").Value = 3023.4615
$ws.Range("N17").Value = -3359.4615

$ws.Range("H62").Value = 2178.724
$ws.Range("I62").Value = 1714.7333
$ws.Range("J62").Value = 2675.8572
$ws.Range("K62").Value = 1714.7333
$ws.Range("L62").Value = 2675.8572
$ws.Range("M62").Value = -1090.7333
$ws.Range("N62").Value = -3923.8572

$ws.Range("H65").Value = 2178.724
$ws.Range("I65").Value = 1714.7333
$ws.Range("J65").Value = 2675.8572
$ws.Range("K65").Value = 8573.666500000001
$ws.Range("L65").Value = 13379.286
$ws.Range("M65").Value = -5453.666500000001
$ws.Range("N65").Value = -19619.286

$ws.Range("H98").Value = 1249.8422
$ws.Range("I98").Value = 1009.1875
$ws.Range("J98").Value = 2533.3333
$ws.Range("K98").Value = 1009.1875
$ws.Range("L98").Value = 2533.3333
$ws.Range("M98").Value = 488.8125
$ws.Range("N98").Value = -5529.3333

$ws.Range("H107").Value = 413.35715
$ws.Range("I107").Value = 202.25
$ws.Range("J107").Value = 1680
$ws.Range("K107").Value = 202.25
$ws.Range("L107").Value = 1680
$ws.Range("M107").Value = 1717.75
$ws.Range("N107").Value = -5520

$ws.Range("H116").Value = 1669.6154
$ws.Range("I116").Value = 1500.625
$ws.Range("J116").Value = 1940
$ws.Range("K116").Value = 1500.625
$ws.Range("L116").Value = 1940
$ws.Range("M116").Value = 1941.375
$ws.Range("N116").Value = -8824

$ws.Range("H122").Value = 1249.8422
$ws.Range("I122").Value = 1009.1875
$ws.Range("J122").Value = 2533.3333
$ws.Range("K122").Value = 3027.5625
$ws.Range("L122").Value = 7599.999899999999
$ws.Range("M122").Value = -577.5625
$ws.Range("N122").Value = -12499.9999

$ws.Range("H132").Value = 702028.8
$ws.Range("I132").Value = 2050.9795
$ws.Range("J132").Value = 2335310.5
$ws.Range("K132").Value = 6152.9385
$ws.Range("L132").Value = 7005931.5
$ws.Range("M132").Value = -3622.9385
$ws.Range("N132").Value = -7010991.5

$ws.Range("H138").Value = 1854094.5
$ws.Range("I138").Value = 1162.0526
$ws.Range("J138").Value = 5054614.5
$ws.Range("K138").Value = 3486.1578
$ws.Range("L138").Value = 15163843.5
$ws.Range("M138").Value = 1653.8422
$ws.Range("N138").Value = -15174123.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1016
$ws.Range("I4").Value = 645
$ws.Range("K4").Value = 645
$ws.Range("M4").Value = -529

$ws.Range("H32").Value = 6974.16
$ws.Range("I32").Value = 3819.6836
$ws.Range("J32").Value = 18841
$ws.Range("K32").Value = 3819.6836
$ws.Range("L32").Value = 18841
$ws.Range("M32").Value = -3532.6836
$ws.Range("N32").Value = -19415

$ws.Range("H107").Value = 15296
$ws.Range("J107").Value = 15296
$ws.Range("L107").Value = 15296
$ws.Range("N107").Value = -22976

$ws.Range("H110").Value = 286455.25
$ws.Range("I110").Value = 357655
$ws.Range("J110").Value = 1656.2858
$ws.Range("K110").Value = 357655
$ws.Range("L110").Value = 1656.2858
$ws.Range("M110").Value = -355610
$ws.Range("N110").Value = -5746.2858

$ws.Range("H122").Value = 2647294
$ws.Range("I122").Value = 1799.4117
$ws.Range("J122").Value = 13890646
$ws.Range("K122").Value = 5398.2351
$ws.Range("L122").Value = 41671938
$ws.Range("M122").Value = -2948.2351
$ws.Range("N122").Value = -41676838

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12898.952
$ws.Range("I86").Value = 18254.143
$ws.Range("J86").Value = 2188.5715
$ws.Range("K86").Value = 18254.143
$ws.Range("L86").Value = 2188.5715
$ws.Range("M86").Value = -17131.143
$ws.Range("N86").Value = -4434.5715

$ws.Range("H89").Value = 12898.952
$ws.Range("I89").Value = 18254.143
$ws.Range("J89").Value = 2188.5715
$ws.Range("K89").Value = 91270.715
$ws.Range("L89").Value = 10942.8575
$ws.Range("M89").Value = -85654.715
$ws.Range("N89").Value = -22174.8575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 21740914
$ws.Range("I58").Value = 29413538
$ws.Range("J58").Value = 1812.75
$ws.Range("K58").Value = 29413538
$ws.Range("L58").Value = 1812.75
$ws.Range("M58").Value = -29413335
$ws.Range("N58").Value = -2218.75

$ws.Range("H99").Value = 11200
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 11200
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").Value = 11200
$ws.Range("N99").Value = -14196

$ws.Range("H126").Value = 11200
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 11200
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").Value = 33600
$ws.Range("N126").Value = -38540

$ws.Range("H134").Value = 17316.373
$ws.Range("I134").Value = 1185.68
$ws.Range("J134").Value = 64759.59
$ws.Range("K134").Value = 3557.04
$ws.Range("L134").Value = 194278.77
$ws.Range("M134").Value = -1022.04
$ws.Range("N134").Value = -199348.77

$ws.Range("H136").Value = 21740914
$ws.Range("I136").Value = 29413538
$ws.Range("J136").Value = 1812.75
$ws.Range("K136").Value = 88240614
$ws.Range("L136").Value = 5438.25
$ws.Range("M136").Value = -88238064
$ws.Range("N136").Value = -10538.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 451.17392
$ws.Range("I113").Value = 338.25
$ws.Range("J113").Value = 574.36365
$ws.Range("K113").Value = 1014.75
$ws.Range("L113").Value = 1723.09095
$ws.Range("M113").Value = 1155.25
$ws.Range("N113").Value = -6063.09095

$ws.Range("H131").Value = 1259.6072
$ws.Range("I131").Value = 363.33334
$ws.Range("J131").Value = 1504.0454
$ws.Range("K131").Value = 1090.00002
$ws.Range("L131").Value = 4512.1362
$ws.Range("M131").Value = 3949.99998
$ws.Range("N131").Value = -14592.1362

$ws.Range("H136").Value = 3136.6667
$ws.Range("I136").Value = 3050
$ws.Range("K136").Value = 9150
$ws.Range("M136").Value = -4050

$ws.Range("H137").Value = 41360
$ws.Range("I137").Value = 1133.3334
$ws.Range("J137").Value = 51416.668
$ws.Range("K137").Value = 3400.0002
$ws.Range("L137").Value = 154250.004
$ws.Range("M137").Value = 1699.9998
$ws.Range("N137").Value = -164450.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 31900
$ws.Range("J59").Value = 31900
$ws.Range("L59").Value = 31900
$ws.Range("N59").Value = -33208

$ws.Range("H132").Value = 23424.39
$ws.Range("I132").Value = 1248.7812
$ws.Range("K132").Value = 3746.3436
$ws.Range("M132").Value = -1216.3436

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 36190.555
$ws.Range("J135").Value = 36190.555
$ws.Range("L135").Value = 36190.555
$ws.Range("N135").Value = -46330.555
